$q3rows = @(
    @("001955","中欧养老产业混合A","22.51","92.70","9.88","2.2240",3),
    @("010429","中欧睿见混合A","18.74","91.22","9.84","1.8440",4),
    @("011710","中欧睿泽混合A","8.33","90.88","9.90","0.8247",3),
    @("166027","中欧创业板两年定期开放混合A","8.11","99.40","9.36","0.7591",4),
    @("004616","中欧电子信息产业沪港深股票A","5.01","92.97","6.65","0.3332",1),
    @("012778","中欧养老产业混合C","2.80","92.70","9.88","0.2766",3),
    @("009791","中欧创业板两年定期开放混合C","2.81","99.40","9.36","0.2630",4),
    @("005763","中欧电子信息产业沪港深股票C","3.88","92.97","6.65","0.2580",1),
    @("003713","英大睿盛灵活配置混合A","2.83","93.65","6.18","0.1749",8),
    @("003714","英大睿盛灵活配置混合C","2.19","93.65","6.18","0.1353",8),
    @("014339","长江智能制造混合A","2.42","75.36","4.13","0.0999",2),
    @("011711","中欧睿泽混合C","0.86","90.88","9.90","0.0851",3),
    @("014155","国泰君安中证500指数增强A","6.64","92.15","1.13","0.0750",10),
    @("014156","国泰君安中证500指数增强C","4.02","92.15","1.13","0.0454",10),
    @("001607","英大策略优选混合A","0.57","91.98","4.56","0.0260",10),
    @("012522","英大稳固增强核心一年持有混合C","1.24","27.71","1.41","0.0175",8),
    @("003447","英大睿鑫灵活配置混合C","0.21","92.71","7.30","0.0153",10),
    @("161038","富国新兴成长量化精选混合（LOF）A","0.81","91.98","1.33","0.0108",9),
    @("012521","英大稳固增强核心一年持有混合A","0.75","27.71","1.41","0.0106",8),
    @("015481","中欧睿见混合C","0.08","91.22","9.84","0.0079",4),
    @("003446","英大睿鑫灵活配置混合A","0.07","92.71","7.30","0.0051",10),
    @("014340","长江智能制造混合C","0.11","75.36","4.13","0.0045",2),
    @("001608","英大策略优选混合C","0.02","91.98","4.56","0.0009",10),
    @("014171","富国新兴成长量化精选混合（LOF）C","0.00","91.98","1.33","0",9),
)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: a new 2022-Q3 row is inserted right after the
#    header, pushing every existing quarter's row down by one.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()

$ws1.Cells.Item(2, 1).Value = 0
$ws1.Cells.Item(3, 1).Copy()
$ws1.Cells.Item(2, 1).PasteSpecial(-4122)

$ws1.Cells.Item(2, 2).Value = "2022-Q3"
$ws1.Cells.Item(2, 3).Value = 24
$ws1.Cells.Item(2, 4).Value = 7.5
$ws1.Cells.Item(2, 2).Style = "Normal"
$ws1.Cells.Item(2, 3).Style = "Normal"
$ws1.Cells.Item(2, 4).Style = "Normal"

for ($r = 3; $r -le 9; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2) New worksheet "2022-Q3" inserted right before the current "2022-Q2"
#    tab (i.e. right after "总计"), holding the per-fund breakdown.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($ws2)
$newSheet.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($j = 0; $j -lt $headers.Count; $j++) {
    $newSheet.Cells.Item(1, $j + 2).Value = $headers[$j]
}

for ($i = 0; $i -lt $q3rows.Count; $i++) {
    $r = $i + 2
    $row = $q3rows[$i]
    $newSheet.Cells.Item($r, 1).Value = $i

    for ($j = 0; $j -lt 5; $j++) {
        $cell = $newSheet.Cells.Item($r, $j + 2)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$j]
    }

    $gcell = $newSheet.Cells.Item($r, 7)
    if ($row[5] -eq "0") {
        $gcell.Value = 0
    } else {
        $gcell.NumberFormat = "@"
        $gcell.Value = $row[5]
    }

    $newSheet.Cells.Item($r, 8).Value = $row[6]
}
